{"js": "// SVM Project Idea Update\n// 1) Add a hanging indent (left=0.5in, hanging=0.5in) to the first\n//    paragraph (\"Machine Learning Project Ideas:\").\n// 2) Replace the \"Vehicle detection and tracking for Autonomous\n//    Vehicles with Support Vector Machines (SVM): ...\" bullet with the\n//    new \"Music Genre Classification ... with Support Vector Machines\n//    (SVM): Build a music recommendation system ...\" text, split across\n//    four runs exactly as the target content does.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- 1) First paragraph gets a hanging indent -----------------------\nconst firstParagraph = paragraphs.items[0];\nfirstParagraph.leftIndent = 36;     // 720 twips = 36 pt = 0.5\"\nfirstParagraph.firstLineIndent = -36; // hanging indent of 0.5\"\nawait context.sync();\n\n// --- 2) Replace the SVM bullet's text/runs ---------------------------\nconst oldText =\n  \"Vehicle detection and tracking for Autonomous Vehicles with Support \" +\n  \"Vector Machines (SVM): Build a model using Support Vector Machines \" +\n  \"(SVM) to diagnose diseases based on medical data and patient symptoms.\";\n\nconst searchResults = body.search(oldText, { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find the target SVM project bullet text.\");\n}\n\nconst targetRange = searchResults.items[0];\n\n// Collapse the matched range to an empty (deleted) range so we can\n// rebuild its contents as four discrete runs.\nconst collapsedRange = targetRange.insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Build the replacement as an OOXML fragment so each <w:t> becomes its\n// own <w:r> run, matching the target structure exactly (rather than\n// collapsing into a single merged run, which insertText would do).\nconst runsOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>' +\n  '<w:r><w:t>Music Genre Classification</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">with Support Vector Machines (SVM): </w:t></w:r>' +\n  '<w:r><w:t>Build a music recommendation system that classifies songs into genres based on audio features, enhancing music streaming platforms.</w:t></w:r>' +\n  '</w:p></w:body></w:document>' +\n  '</pkg:xmlData></pkg:part>' +\n  '</pkg:package>';\n\ncollapsedRange.insertOoxml(runsOoxml, Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "# SVM Project Idea Update\n# 1) Add a hanging indent (left=0.5in, hanging=0.5in) to the first\n#    paragraph (\"Machine Learning Project Ideas:\").\n# 2) Replace the \"Vehicle detection and tracking for Autonomous\n#    Vehicles with Support Vector Machines (SVM): ...\" bullet with the\n#    new \"Music Genre Classification ... with Support Vector Machines\n#    (SVM): Build a music recommendation system ...\" text, split across\n#    four runs exactly as the target content does.\n\n$d = $word.ActiveDocument\n\n# --- 1) First paragraph gets a hanging indent -------------------------\n$firstPara = $d.Paragraphs(1)\n$firstPara.Range.ParagraphFormat.LeftIndent = 36       # 720 twips = 36 pt = 0.5\"\n$firstPara.Range.ParagraphFormat.FirstLineIndent = -36 # hanging indent of 0.5\"\n\n# --- 2) Replace the SVM bullet's text/runs -----------------------------\n$oldText = \"Vehicle detection and tracking for Autonomous Vehicles with Support Vector Machines (SVM): Build a model using Support Vector Machines (SVM) to diagnose diseases based on medical data and patient symptoms.\"\n\n$searchRange = $d.Content\n$searchRange.Find.Text = $oldText\n$found = $searchRange.Find.Execute()\n\nif (-not $found) {\n    throw \"Could not find the target SVM project bullet text.\"\n}\n\n# Work against the whole enclosing paragraph so its pPr (style + list\n# numbering) is preserved, and so the new runs end up in place of the\n# old ones (with the trailing line-break run kept at the end).\n$targetParagraph = $searchRange.Paragraphs(1)\n$paraRange = $targetParagraph.Range\n\n$ooxml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr></w:pPr><w:r><w:t>Music Genre Classification</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t xml:space=\"preserve\">with Support Vector Machines (SVM): </w:t></w:r><w:r><w:t>Build a music recommendation system that classifies songs into genres based on audio features, enhancing music streaming platforms.</w:t></w:r><w:r><w:br/></w:r></w:p></w:body></w:document>\n</pkg:xmlData></pkg:part>\n</pkg:package>\n'@\n\n$paraRange.InsertXML($ooxml)\n"}
